$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header "Save" in H1, copy style from G1 (bold, centered, bordered)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set values for rows 2-5 in column H
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
